$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "auto usate"

$ws.Range("A3").Value = "Nome Colonna : "
$ws.Range("A4").Value = "Tipo di Dato :"
$ws.Range("A5").Value = "Chiave Primaria :"

$ws.Range("B3").Value = "modello"
$ws.Range("C3").Value = "porte"
$ws.Range("D3").Value = "numero_Sedie"
$ws.Range("E3").Value = "targa"
$ws.Range("F3").Value = "numero_proprietari_precendti"
$ws.Range("G3").Value = "cilindrata"
$ws.Range("H3").Value = "uso_commerciale"
$ws.Range("I3").Value = "motore_originale"

$ws.Range("B3:Q3").Font.Bold = $true

$ws.Columns.Item(1).ColumnWidth = 15.5
$ws.Columns.Item(4).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 27.333333333333332
$ws.Columns.Item(8).ColumnWidth = 15.333333333333334
$ws.Columns.Item(9).ColumnWidth = 15.5

$ws.PageSetup.Orientation = 1

$ws.Range("L9").Select()
